$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (ID Competição) for rows 2-71 was incorrectly truncated to 37;
# restore the dropped leading digit so the value reads 237.
$ws.Range("B2:B71").Value = 237
